# Update cfb_weather.xlsx with Timestamp 2024-09-22T16:21:36.405691
$wb = $excel.ActiveWorkbook

# --- 1) Refresh the run "Timestamp" column on the FBS sheet ---
# Every data row (2-52) in column AD shares the same timestamp string;
# rewrite them all to the new capture time.
$wsFbs = $wb.Worksheets.Item("FBS")
$newTimestamp = "2024-09-22T16:21:36.405691"
$lastRow = $wsFbs.Cells.Item($wsFbs.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $wsFbs.Range("AD" + $r).Value = $newTimestamp
}

# --- 2) Update the recalculated weather figures on the "Other" sheet (row 2) ---
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("K2").Value = 79.09999999999999   # temp_fg
$wsOther.Range("L2").Value = 7.6                  # wind_fg
$wsOther.Range("P2").Value = -6.4                 # wind_diff
